$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 139, shifting existing rows 139+ down to 140+.
$ws.Rows.Item(139).Insert()

# Populate the new row 139 with a new weekly entry (same record as the row
# that follows it, but one day later and with a different Volumen value).
$ws.Range("A139").Value = 4
$ws.Range("B139").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C139").Value = "Los Lagos"
$ws.Range("D139").Value = 44587
$ws.Range("E139").Value = 10
$ws.Range("F139").Value = 100112023
$ws.Range("G139").Value = "Brócoli"
$ws.Range("H139").Value = "Sin especificar"
$ws.Range("I139").Value = "Primera"
$ws.Range("J139").Value = 100
$ws.Range("K139").Value = 1500
$ws.Range("L139").Value = 1500
$ws.Range("M139").Value = 1500
$ws.Range("N139").Value = "$/unidad"
$ws.Range("O139").Value = "Región Metropolitana"
$ws.Range("P139").Value = 1500
$ws.Range("Q139").Value = 1
$ws.Range("R139").Value = "Hortaliza"
